$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 2).Value = 0.0090293453724605
$ws.Cells.Item(2, 3).Value = 0.000752445447705041
$ws.Cells.Item(2, 4).Value = 0.00526711813393529
$ws.Cells.Item(2, 5).Value = 0.00752445447705041
$ws.Cells.Item(2, 6).Value = 0.000752445447705041
$ws.Cells.Item(2, 7).Value = 0.00526711813393529
$ws.Cells.Item(2, 8).Value = 0.936794582392777
$ws.Cells.Item(2, 9).Value = 0.00526711813393529
$ws.Cells.Item(2, 10).Value = 0.00376222723852521
$ws.Cells.Item(2, 11).Value = 0.0090293453724605
$ws.Cells.Item(2, 12).Value = 0.867569601203913
$ws.Cells.Item(2, 13).Value = 0.151241534988713
$ws.Cells.Item(2, 14).Value = 0.0353649360421369
$ws.Cells.Item(2, 15).Value = 0.0112866817155756
$ws.Cells.Item(2, 16).Value = 0.00376222723852521
$ws.Cells.Item(2, 17).Value = 0.0090293453724605
$ws.Cells.Item(2, 18).Value = 0.0127915726109857
$ws.Cells.Item(2, 19).Value = 0.09104589917231
$ws.Cells.Item(2, 20).Value = 0.0737396538750941
$ws.Cells.Item(2, 21).Value = 0.018058690744921
$ws.Cells.Item(2, 22).Value = 0.00376222723852521
$ws.Cells.Item(2, 23).Value = 0.0158013544018059
$ws.Cells.Item(2, 24).Value = 0.0090293453724605
$ws.Cells.Item(3, 2).Value = 0.973664409330323
$ws.Cells.Item(3, 3).Value = 0.000752445447705041
$ws.Cells.Item(3, 4).Value = 0.00225733634311512
$ws.Cells.Item(3, 5).Value = 0.00677200902934537
$ws.Cells.Item(3, 6).Value = 0.99849510910459
$ws.Cells.Item(3, 7).Value = 0.99398043641836
$ws.Cells.Item(3, 8).Value = 0
$ws.Cells.Item(3, 9).Value = 0
$ws.Cells.Item(3, 10).Value = 0
$ws.Cells.Item(3, 11).Value = 0
$ws.Cells.Item(3, 12).Value = 0.00376222723852521
$ws.Cells.Item(3, 13).Value = 0.00376222723852521
$ws.Cells.Item(3, 14).Value = 0.0135440180586907
$ws.Cells.Item(3, 15).Value = 0.528216704288939
$ws.Cells.Item(3, 16).Value = 0.124153498871332
$ws.Cells.Item(3, 17).Value = 0.0684725357411588
$ws.Cells.Item(3, 18).Value = 0.00827689992475546
$ws.Cells.Item(3, 19).Value = 0.00752445447705041
$ws.Cells.Item(3, 20).Value = 0.735891647855531
$ws.Cells.Item(3, 21).Value = 0.929270127915726
$ws.Cells.Item(3, 22).Value = 0.188863807373965
$ws.Cells.Item(3, 23).Value = 0.0353649360421369
$ws.Cells.Item(3, 24).Value = 0.00677200902934537
$ws.Cells.Item(4, 2).Value = 0.00451467268623025
$ws.Cells.Item(4, 3).Value = 0.0120391271632807
$ws.Cells.Item(4, 4).Value = 0.0112866817155756
$ws.Cells.Item(4, 5).Value = 0.0112866817155756
$ws.Cells.Item(4, 6).Value = 0
$ws.Cells.Item(4, 7).Value = 0.000752445447705041
$ws.Cells.Item(4, 8).Value = 0.0564334085778781
$ws.Cells.Item(4, 9).Value = 0.99097065462754
$ws.Cells.Item(4, 10).Value = 0.996237772761475
$ws.Cells.Item(4, 11).Value = 0.990218209179834
$ws.Cells.Item(4, 12).Value = 0.126410835214447
$ws.Cells.Item(4, 13).Value = 0.840481565086531
$ws.Cells.Item(4, 14).Value = 0.947328818660647
$ws.Cells.Item(4, 15).Value = 0.18660647103085
$ws.Cells.Item(4, 16).Value = 0.838976674191121
$ws.Cells.Item(4, 17).Value = 0.161775771256584
$ws.Cells.Item(4, 18).Value = 0.826937547027841
$ws.Cells.Item(4, 19).Value = 0.862302483069977
$ws.Cells.Item(4, 20).Value = 0.155756207674944
$ws.Cells.Item(4, 21).Value = 0.036117381489842
$ws.Cells.Item(4, 22).Value = 0.0285929270127916
$ws.Cells.Item(4, 23).Value = 0.798344620015049
$ws.Cells.Item(4, 24).Value = 0.948833709556057
$ws.Cells.Item(5, 2).Value = 0.0127915726109857
$ws.Cells.Item(5, 3).Value = 0.986455981941309
$ws.Cells.Item(5, 4).Value = 0.981188863807374
$ws.Cells.Item(5, 5).Value = 0.974416854778029
$ws.Cells.Item(5, 6).Value = 0.000752445447705041
$ws.Cells.Item(5, 7).Value = 0
$ws.Cells.Item(5, 8).Value = 0.00677200902934537
$ws.Cells.Item(5, 9).Value = 0.00376222723852521
$ws.Cells.Item(5, 10).Value = 0
$ws.Cells.Item(5, 11).Value = 0.000752445447705041
$ws.Cells.Item(5, 12).Value = 0.00225733634311512
$ws.Cells.Item(5, 13).Value = 0.00451467268623025
$ws.Cells.Item(5, 14).Value = 0.00376222723852521
$ws.Cells.Item(5, 15).Value = 0.27313769751693
$ws.Cells.Item(5, 16).Value = 0.0331075996990218
$ws.Cells.Item(5, 17).Value = 0.760722347629797
$ws.Cells.Item(5, 18).Value = 0.151993980436418
$ws.Cells.Item(5, 19).Value = 0.0391271632806622
$ws.Cells.Item(5, 20).Value = 0.0338600451467269
$ws.Cells.Item(5, 21).Value = 0.0158013544018059
$ws.Cells.Item(5, 22).Value = 0.778781038374718
$ws.Cells.Item(5, 23).Value = 0.150489089541008
$ws.Cells.Item(5, 24).Value = 0.0353649360421369
